$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row of AWS / Databricks lab credentials (S.No 24) below the
# existing 23 rows of data (header is row 1, data rows 2-24).
$row = 25

$ws.Cells.Item($row, 1).Value = 24
$ws.Cells.Item($row, 2).Value = "clouduser"
$ws.Cells.Item($row, 3).Value = "YwcY5IkyF2AlxG5M"
$ws.Cells.Item($row, 4).Value = "https://883601963544.signin.aws.amazon.com/console?region=us-west-2"
$ws.Cells.Item($row, 5).Value = "us-west-2"
$ws.Cells.Item($row, 6).Value = "https://accounts.cloud.databricks.com/"
$ws.Cells.Item($row, 7).Value = "akhila.gudla6@gmail.com"

# New row's look: white fill, medium black box border, wrap text, vertically
# centered, small 6pt Calibri font for the plain cells.
$rowRange = $ws.Range("A25:G25")
$rowRange.Interior.Color = 16777215
$rowRange.Borders.Weight = -4138
$rowRange.Borders.Color = 0
$rowRange.WrapText = $true
$rowRange.VerticalAlignment = -4108
$rowRange.HorizontalAlignment = -4131
$rowRange.Font.Name = "Calibri"
$rowRange.Font.Size = 6
$rowRange.Font.Color = 0

# Hyperlinks for the login URL, console URL and e-mail columns, matching the
# hyperlinks already present on the earlier rows.
$ws.Hyperlinks.Add($ws.Range("D25"), "https://883601963544.signin.aws.amazon.com/console?region=us-west-2", "", "https://883601963544.signin.aws.amazon.com/console?region=us-west-2") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F25"), "https://accounts.cloud.databricks.com/", "", "https://accounts.cloud.databricks.com/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G25"), "mailto:akhila.gudla6@gmail.com", "", "mailto:akhila.gudla6@gmail.com") | Out-Null

# Hyperlinked cells keep the left-aligned wrap-text look, except the login
# URL cell (D25) which is centered.
$ws.Range("D25").HorizontalAlignment = -4108
$ws.Range("F25:G25").HorizontalAlignment = -4131

# Update the used-range dimension and selection to include the new row.
$ws.Range("D25").Select()
